$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-04-04 Thursday", $false, $false, $false, $false, $false, $true, 1, $false, "2024-04-05 Friday", 2)

# Update the division expressions in the table (row -> array of 5 new values)
$table = $d.Tables.Item(1)

$rowMap = @{
    1  = @("96÷4=", "89÷9=", "17÷4=", "62÷6=", "68÷9=")
    5  = @("84÷7=", "42÷5=", "29÷4=", "18÷9=", "90÷8=")
    9  = @("21÷2=", "33÷3=", "56÷3=", "40÷9=", "13÷8=")
    13 = @("88÷2=", "75÷8=", "87÷6=", "66÷4=", "42÷6=")
    17 = @("64÷5=", "53÷2=", "25÷8=", "45÷3=", "74÷8=")
}

foreach ($rowIndex in $rowMap.Keys) {
    $values = $rowMap[$rowIndex]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $table.Cell($rowIndex, $col)
        $cellRange = $cell.Range
        $cellRange.End = $cellRange.End - 1
        $cellRange.Text = $values[$col - 1]
    }
}
